$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the sheet, shifting existing rows down.
$ws.Rows.Item(1).Insert()

# Populate the new header row.
$ws.Cells.Item(1, 1).Value2 = "Parameter notation"
$ws.Cells.Item(1, 2).Value2 = "Definition"

# Update the view selection / scroll position to match the authored state.
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
$ws.Range("A43").Select() | Out-Null
